$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell "D2" "89.380.97"
Set-TextCell "E2" "  +10.42%  "
Set-TextCell "D3" "3.362.18"
Set-TextCell "E3" "  +7.15%  "
Set-TextCell "E4" "  +0.07%  "
Set-TextCell "D5" "219.45"
Set-TextCell "E5" "  +5.58%  "
Set-TextCell "D6" "647.64"
Set-TextCell "E6" "  +5.15%  "
Set-TextCell "D7" "0.412"
Set-TextCell "E7" "  +46.82%  "
Set-TextCell "D8" "0.999"
Set-TextCell "E8" "  +0.00%  "
Set-TextCell "E9" "  +6.87%  "
Set-TextCell "D10" "3.358.99"
Set-TextCell "E10" "  +7.25%  "
Set-TextCell "D11" "0.613"
Set-TextCell "E11" "  +7.35%  "
Set-TextCell "D12" "0.0000290"
Set-TextCell "E12" "  +15.78%  "
Set-TextCell "D13" "36.48"
Set-TextCell "E13" "  +16.82%  "
Set-TextCell "E14" "  +2.41%  "
Set-TextCell "D15" "3.981.72"
Set-TextCell "E15" "  +7.28%  "
Set-TextCell "D16" "5.55"
Set-TextCell "E16" "  +5.66%  "
Set-TextCell "D17" "89.280.21"
Set-TextCell "E17" "  +10.51%  "
Set-TextCell "D18" "3.358.51"
Set-TextCell "E18" "  +7.28%  "
Set-TextCell "D19" "14.81"
Set-TextCell "E19" "  +7.49%  "
Set-TextCell "D20" "3.13"
Set-TextCell "E20" "  +0.46%  "
Set-TextCell "D21" "9.69"
Set-TextCell "E21" "  +8.84%  "
Set-TextCell "D22" "457.74"
Set-TextCell "E22" "  +6.88%  "
Set-TextCell "E23" "  +9.94%  "
Set-TextCell "D24" "7.51"
Set-TextCell "E24" "  +4.94%  "
Set-TextCell "D25" "5.57"
Set-TextCell "E25" "  +8.18%  "
Set-TextCell "D26" "12.82"
Set-TextCell "E26" "  +19.36%  "
Set-TextCell "D27" "3.524.33"
Set-TextCell "E27" "  +6.82%  "
Set-TextCell "E28" "  +19.00%  "
Set-TextCell "D29" "79.23"
Set-TextCell "E29" "  +4.98%  "
Set-TextCell "D30" "0.200"
Set-TextCell "E30" "  +45.08%  "
Set-TextCell "D31" "0.998"
Set-TextCell "E31" "  -0.02%  "
Set-TextCell "D32" "9.44"
Set-TextCell "E32" "  +6.24%  "
Set-TextCell "D33" "596.01"
Set-TextCell "E33" "  +7.16%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextCell "D34" "0.999"
Set-TextCell "E34" "  +0.06%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell "D35" "1.57"
Set-TextCell "E35" "  +7.30%  "
Set-TextCell "E36" "  +7.24%  "
Set-TextCell "D37" "7.34"
Set-TextCell "E37" "  +22.37%  "
Set-TextCell "E38" "  -4.57%  "
Set-TextCell "D39" "23.53"
Set-TextCell "E39" "  +4.78%  "
Set-TextCell "D40" "0.429"
Set-TextCell "E40" "  +6.42%  "
Set-TextCell "D41" "2.15"
Set-TextCell "E41" "  +7.12%  "
$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextCell "D42" "21.85"
Set-TextCell "E42" "  +5.48%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell "D43" "3.20"
Set-TextCell "E43" "  +6.38%  "
Set-TextCell "D44" "0.998"
Set-TextCell "E44" "  -0.05%  "
Set-TextCell "D45" "1.47"
Set-TextCell "E45" "  +12.21%  "
Set-TextCell "D46" "158.17"
Set-TextCell "E46" "  -0.52%  "
Set-TextCell "E47" "  +0.06%  "
Set-TextCell "D48" "189.96"
Set-TextCell "E48" "  +1.89%  "
Set-TextCell "D49" "46.32"
Set-TextCell "E49" "  +4.16%  "
Set-TextCell "D50" "4.51"
Set-TextCell "E50" "  +8.37%  "
Set-TextCell "D51" "0.667"
Set-TextCell "E51" "  +7.66%  "
